$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.Style = "Normal"
}

Set-TextValue "D2" "298.42"
Set-TextValue "E2" "-2.22%"
Set-TextValue "D3" "31.46"
Set-TextValue "E3" "-1.45%"
Set-TextValue "D4" "5.152"
Set-TextValue "E4" "-2.35%"
Set-TextValue "D5" "0.07318"
Set-TextValue "E5" "-2.56%"
Set-TextValue "D6" "1.866"
Set-TextValue "E6" "26.53%"
Set-TextValue "D7" "7.763"
Set-TextValue "E7" "-1.04%"
Set-TextValue "D8" "3.741"
Set-TextValue "E8" "-0.49%"
Set-TextValue "D9" "0.9255"
Set-TextValue "E9" "1.05%"
Set-TextValue "D10" "0.1676"
Set-TextValue "E10" "-0.99%"
Set-TextValue "D11" "0.07154"
Set-TextValue "E11" "-8.43%"
Set-TextValue "D12" "0.07976"
Set-TextValue "E12" "-0.83%"
Set-TextValue "D13" "0.02998"
Set-TextValue "E13" "-0.58%"
Set-TextValue "D14" "0.09924"
Set-TextValue "E14" "0.32%"
Set-TextValue "D15" "0.001488"
Set-TextValue "E15" "-0.32%"
Set-TextValue "D16" "0.006163"
Set-TextValue "E16" "-2.09%"
Set-TextValue "D17" "3.456"
Set-TextValue "E17" "-0.77%"
Set-TextValue "E18" "-0.53%"
Set-TextValue "E19" "-2.28%"
Set-TextValue "E20" "-1.97%"
Set-TextValue "D21" "4.559"
Set-TextValue "E21" "1.93%"
Set-TextValue "D22" "0.04641"
Set-TextValue "E22" "1.86%"
Set-TextValue "D23" "0.1582"
Set-TextValue "E23" "-3.34%"
Set-TextValue "D24" "0.001216"
Set-TextValue "E24" "0.00%"
Set-TextValue "D25" "0.004735"
Set-TextValue "E25" "6.84%"
Set-TextValue "E26" "-7.06%"
Set-TextValue "D27" "0.0001874"
Set-TextValue "E27" "7.74%"
Set-TextValue "D39" "0.01715"
Set-TextValue "E39" "-2.13%"
Set-TextValue "D40" "0.04475"
Set-TextValue "E40" "-1.16%"
Set-TextValue "D41" "0.007061"
Set-TextValue "E41" "-2.40%"
Set-TextValue "D42" "0.1330"
Set-TextValue "E42" "-1.11%"
Set-TextValue "D43" "0.002148"
Set-TextValue "E43" "-3.93%"
Set-TextValue "E44" "-24.31%"
Set-TextValue "D45" "0.00006230"
Set-TextValue "E45" "0.34%"
Set-TextValue "E46" "-21.35%"
Set-TextValue "D47" "1.920"
Set-TextValue "E47" "171.13%"

Write-Host "Updated 64 cells"
